$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for match "Kz2IMp1S" (Colorado Rapids - LA Galaxy) was removed.
# Deleting row 2 shifts row 3 (rcPS0RFN / Las Vegas Lights - Sacramento
# Republic) up into row 2, which already carries most of the correct data.
$ws.Rows.Item(2).Delete()

# A handful of odds for the now-row-2 match were refreshed; patch just those.
$ws.Range("G2").Value = 2.65
$ws.Range("H2").Value = 2.67
$ws.Range("I2").Value = 2.87
$ws.Range("K2").Value = 1.93
$ws.Range("L2").Value = 3.5
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 5.8
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.91
$ws.Range("W2").Value = 7.4
$ws.Range("AB2").Value = 35
$ws.Range("AC2").Value = 5.8
$ws.Range("AD2").Value = 5.3
$ws.Range("AE2").Value = 13
$ws.Range("AI2").Value = 15
$ws.Range("AL2").Value = 27
$ws.Range("AN2").Value = 4.6
$ws.Range("AP2").Value = 22
$ws.Range("AU2").Value = 6.5
$ws.Range("AW2").Value = 4.85
$ws.Range("AX2").Value = 16.5
$ws.Range("BB2").Value = 300
